$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.388.66'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.63%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.012.49'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.65%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.11%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5138'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.61%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4264'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.22%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08750'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.81%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.54'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.78%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.135'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.76%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.45'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.42%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.014.33'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.72%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.644'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.64%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.466'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.98%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.52%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.13%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001113'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.43%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06550'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.50%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.16%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.200'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.446.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.59%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.39%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.252'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.79%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.253.71'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.32%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.67%  '

$ws.Range('E28').Value = '  -0.12%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.458'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.17%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '131.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.98%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.139'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.68%  '

$ws.Range('E32').Value = '  +1.68%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.086'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.13%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.829'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.49%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.364'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.98%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02540'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.68%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06677'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.99%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.462'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.86%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.45'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.25%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.215'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.12%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2217'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.82%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6667'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.85%  '

$ws.Range('E43').Value = '  +2.31%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.11%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.70'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.60%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6182'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.26%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.205'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.17%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.635'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.53%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.260'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.07%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.70%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.26'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.81%  '
